# Fruta / hortaliza, semanal
# A new weekly data point (row) is inserted at the top of the existing
# Mango / Femacal de La Calera block (row 291), pushing all the
# following rows (291-373) down by one (to 292-374). The new row
# reuses the same Mercado/Region/Producto/Calidad/Volumen/Precio/Unidad
# data that was already present for that slot, only the date and the
# origin differ for the newest entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 291:373 down to 292:374, leaving a blank row 291 behind.
$ws.Rows(291).Insert()

# Populate the newly inserted row 291 with the latest weekly record.
$ws.Cells.Item(291, 1).Value = 3
$ws.Cells.Item(291, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(291, 3).Value = "Coquimbo"
$ws.Cells.Item(291, 4).Value = 44642
$ws.Cells.Item(291, 5).Value = 5
$ws.Cells.Item(291, 6).Value = "Fruta"
$ws.Cells.Item(291, 7).Value = 100108
$ws.Cells.Item(291, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(291, 9).Value = 100108002
$ws.Cells.Item(291, 10).Value = "Mango"
$ws.Cells.Item(291, 11).Value = "Sin especificar"
$ws.Cells.Item(291, 12).Value = "Primera"
$ws.Cells.Item(291, 13).Value = 228
$ws.Cells.Item(291, 14).Value = 9000
$ws.Cells.Item(291, 15).Value = 9000
$ws.Cells.Item(291, 16).Value = 9000
$ws.Cells.Item(291, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(291, 18).Value = "Perú"
$ws.Cells.Item(291, 19).Value = 2250
$ws.Cells.Item(291, 20).Value = 4

# Match the date cell's number format to the rest of the date column.
$ws.Cells.Item(291, 4).NumberFormat = $ws.Cells.Item(292, 4).NumberFormat
